# QuickRoute.docx edit — wording / proofreading pass
# ("En lille smule rettelser.. o.O")
#
# The underlying body text is rewritten: the opening sentence is
# restructured, a couple of word choices are corrected ("løbt" ->
# "løbet", "proffessionel" -> "professionel"), and a passive/active
# voice tweak is made in the Google Earth sentence. We apply each
# change as a targeted Find/Replace over the whole story so the
# surrounding run formatting (bold heading, da-DK language tag, etc.)
# is preserved untouched.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $range.Find.ClearFormatting()
    $range.Find.Replacement.ClearFormatting()
    [void]$range.Find.Execute(
        $find,      # FindText
        $false,     # MatchCase
        $false,     # MatchWholeWord
        $false,     # MatchWildcards
        $false,     # MatchSoundsLike
        $false,     # MatchAllWordForms
        $true,      # Forward
        1,          # Wrap (wdFindContinue)
        $false,     # Format
        $replace,   # ReplaceWith
        2           # Replace (wdReplaceAll)
    )
}

# 1. Opening sentence of the body paragraph is restructured.
Replace-Text `
    "Som Claus beskrev i vores interview, er QuickRoute en eksisterende løsning på" `
    "Igennem gruppens interview, beskrev Claus, at QuickRoute er en eksisterende løsning på"

# 2. "løbt" -> "løbet"
Replace-Text `
    "samlet distance løbt til det punkt" `
    "samlet distance løbet til det punkt"

# 3. Google Earth sentence: active -> passive voice tweak.
Replace-Text `
    "Til at toppe alt dette af, kan man efterfølgende integrere det på Google Earth, så man kan se en 3D model" `
    "Til at toppe alt dette af, kan det efterfølgende integreres på Google Earth, så der kan ses en 3D model"

# 4. Typo fix: "proffessionel" -> "professionel"
Replace-Text `
    "proffessionel" `
    "professionel"
